$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(15).Insert()

$ws.Range("A15").Value = "PROFIL"
$ws.Range("B15").Value = "O"
$ws.Range("C15").Value = "Ocynkownia"
$ws.Range("D15").Value = "Ocynk"
$ws.Range("E15").Value = "Centrum kompletacji"
$ws.Range("F15").Value = "Kompletacja"
$ws.Range("G15").Value = "Montaż"
$ws.Range("H15").Value = "Biały_Montaż"

$ws.Range("B15").Select()
